$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.81409999999998
$ws.Range("A6").Value = -22.78780000000001
$ws.Range("A7").Value = -20.27869999999997
$ws.Range("A8").Value = -22.36700000000002
$ws.Range("A16").Value = -21.96800000000001
$ws.Range("A20").Value = -19.9414
$ws.Range("A21").Value = -20.19049999999998
$ws.Range("A28").Value = -21.87279999999999
$ws.Range("A29").Value = -21.40679999999996
$ws.Range("A30").Value = -21.44580000000001
$ws.Range("A32").Value = -21.12609999999999
$ws.Range("A40").Value = -20.12590000000001
$ws.Range("A46").Value = -21.82980000000001
$ws.Range("A51").Value = -21.54869999999998
$ws.Range("A52").Value = -22.12779999999999
$ws.Range("A57").Value = -22.68640000000003
$ws.Range("A59").Value = -22.22420000000001
$ws.Range("A62").Value = -22.15020000000001
$ws.Range("A66").Value = -21.4324
$ws.Range("A73").Value = -20.2978
$ws.Range("A74").Value = -21.63809999999998
$ws.Range("A77").Value = -20.11259999999999
$ws.Range("A92").Value = -21.6292
$ws.Range("A100").Value = -22.08660000000001
